# Update TPM-derived NATMI values for Wnt5a-Lrp5 LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.07063700000000001
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("M2").Value = 15.01856033333333
$ws.Range("N2").Value = 45.055681
$ws.Range("O2").Value = 0.4908713633047416
$ws.Range("P2").Value = 0.4908713633047417
$ws.Range("Q2").Value = 0.3536220154218889
$ws.Range("R2").Value = 3.182598138797
$ws.Range("S2").Value = 0.00138207551396958
$ws.Range("T2").Value = 0.00138207551396958
$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.07063700000000001
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("O3").Value = 0.3099803572711625
$ws.Range("P3").Value = 0.3099803572711625
$ws.Range("Q3").Value = 0.2233087665604445
$ws.Range("R3").Value = 2.009778899044
$ws.Range("S3").Value = 0.0008727668664795332
$ws.Range("T3").Value = 0.0008727668664795332
$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.07063700000000001
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("O4").Value = 0.1991482794240958
$ws.Range("P4").Value = 0.1991482794240958
$ws.Range("Q4").Value = 0.143465724836
$ws.Range("R4").Value = 1.291191523524
$ws.Range("S4").Value = 0.0005607130120368055
$ws.Range("T4").Value = 0.0005607130120368055
$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("M5").Value = 15.01856033333333
$ws.Range("N5").Value = 45.055681
$ws.Range("O5").Value = 0.4908713633047416
$ws.Range("P5").Value = 0.4908713633047417
$ws.Range("Q5").Value = 123.9436978938339
$ws.Range("R5").Value = 1115.493281044505
$ws.Range("S5").Value = 0.4844142686239198
$ws.Range("T5").Value = 0.4844142686239198
$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("O6").Value = 0.3099803572711625
$ws.Range("P6").Value = 0.3099803572711625
$ws.Range("S6").Value = 0.3059027665504093
$ws.Range("T6").Value = 0.3059027665504093
$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("O7").Value = 0.1991482794240958
$ws.Range("P7").Value = 0.1991482794240958
$ws.Range("S7").Value = 0.1965286128639876
$ws.Range("T7").Value = 0.1965286128639876
$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("M8").Value = 15.01856033333333
$ws.Range("N8").Value = 45.055681
$ws.Range("O8").Value = 0.4908713633047416
$ws.Range("P8").Value = 0.4908713633047417
$ws.Range("Q8").Value = 1.298509732606778
$ws.Range("R8").Value = 11.686587593461
$ws.Range("S8").Value = 0.005075019166852266
$ws.Range("T8").Value = 0.005075019166852266
$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("O9").Value = 0.3099803572711625
$ws.Range("P9").Value = 0.3099803572711625
$ws.Range("S9").Value = 0.003204823854273649
$ws.Range("T9").Value = 0.00320482385427365
$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("O10").Value = 0.1991482794240958
$ws.Range("P10").Value = 0.1991482794240958
$ws.Range("S10").Value = 0.002058953548071388
$ws.Range("T10").Value = 0.002058953548071388
